$p = $ppt.ActivePresentation

# Slides 6, 7 and 8 each have an empty title placeholder (shape 1) whose
# paragraph currently holds only an <a:endParaRPr sz="2400"/> and no run.
# Add a single "." run at 24pt to that paragraph on each of those slides,
# while keeping the existing end-of-paragraph run properties intact.
foreach ($idx in 6, 7, 8) {
    $s = $p.Slides.Item($idx)
    $titleShape = $s.Shapes.Title
    $tr = $titleShape.TextFrame.TextRange

    # Setting the font size on the still-empty text range first preserves
    # the paragraph's trailing end-of-paragraph run properties; only then
    # assign the text itself, which becomes the new run.
    $tr.Font.Size = 24
    $tr.Text = "."
}
